# Minor fix in TSP.
# Update the "Fitness" values (column C) for rows 2-12 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value  = 3927
$ws.Range("C3").Value  = 4109
$ws.Range("C4").Value  = 4109
$ws.Range("C5").Value  = 4119
$ws.Range("C6").Value  = 4119
$ws.Range("C7").Value  = 4119
$ws.Range("C8").Value  = 4666
$ws.Range("C9").Value  = 4666
$ws.Range("C10").Value = 4666
$ws.Range("C11").Value = 4790
$ws.Range("C12").Value = 4928
